$d = $word.ActiveDocument

# 1. Activation date change
$d.Content.Find.Execute("Ativação: 01/01/2021", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2024", 2) | Out-Null

# 2. Portuguese objectives paragraph
$old2 = "Apresentar  os conceitos básicos da Ciência Econômica, capacitando-os a compreender os principais conceitos micro e macroeconômicos e a interpretar o discurso e a prática da economia, orientados pelo seu próprio senso crítico."
$new2 = "Apresentar os conceitos básicos da Ciência Econômica, capacitando os alunos a compreender os principais conceitos micro e macroeconômicos, e a analisar o discurso e a prática da economia, orientados pelo seu próprio senso crítico."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3. English objectives paragraph
$old3 = "Introduce the students of Production Engineering to the basic concepts of Economic Science, enabling students to understand the main micro and macroeconomic concepts and to interpret the discourse and practice of economics, guided by their own critical sense."
$new3 = "Introduce the basic concepts of Economic Science, enabling the students to understand the main micro and macroeconomic concepts and to analyze the discourse and practice of economics, guided by their critical sense."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# 4. Add a new docente after Herlandí's entry
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("11079086 - Herlandí de Souza Andrade")) {
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Collapse(0)  # collapse to end (wdCollapseEnd = 0), excluding paragraph mark
        $boundary = $r.Start
        $r.InsertAfter("3295113 - José Eduardo Holler Branco")
        $rb = $d.Range($boundary, $boundary)
        $rb.InsertBreak(6)  # wdLineBreak
        break
    }
}

# 5. Método text
$old5 = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$new5 = "Provas, trabalhos em grupo, exercícios individuais e seminários."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# 6. Critério text
$old6 = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas."
$new6 = "Média das atividades avaliativas."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# 7. Bibliography paragraph
$old7 = "VASCONCELLOS, M. A. S.; GARCIA, M. E. Fundamentos de Economia. 6 ed. São Paulo: Saraiva, 2018.GREMAUD, A. P. Introdução à Economia. São Paulo: Atlas, 2017.ROSSETTI, J. P. Introdução à Economia - Livro Texto. São Paulo: Atlas, 2016.VASCONCELLOS, M. A. S. ECONOMIA: Micro e Macro. São Paulo: Atlas, 2015.ALBERGONI, L. INTRODUÇÃO À ECONOMIA: Aplicações no Cotidiano. São Paulo: Atlas, 2015.GREMAUD, A. P.; VASCONCELLOS, M. A. S.; TONETO JÚNIOR, R. Economia Brasileira Contemporânea. 8 ed. São Paulo: Atlas, 2017.MÉNARD, C.; SAES, M. S. M.; SILVA, V. L. S.; RAYNAUD, E. Economia das Organizações: Formas Plurais e Desafios. São Paulo: Atlas, 2014.BACHA et Al. Estado da Economia Mundial - Desafios e Respostas - Seminário em Homenagem a Pedro Malan. São Paulo: LTC, 2015.BACHA , Edmar. Introdução à Macroeconomia: Uma perspectiva brasileira. Rio de Janeiro: Campus,1987.BEGG, D.; DORNBUSCH, R.; FISCHER, S. Introdução A Economia. Rio de Janeiro: Campus, 2003. FURTADO, C. Formação econômica do Brasil. São Paulo: Companhia Editora Nacional, 2003.GRAMAUD, A. P. et alli. Manual de economia. São Paulo. Saraiva. 2004.GRAMAUD, A. P. et alli. Economia Brasileira Contemporânea. 6.ed. São Paulo. Atlas, 2006.HUNT, E. K.; SHERMAN, H. J. História do Pensamento Econômico. Petrópolis : Vozes, 1997.MANKIW, N.G. Introdução à economia. São Paulo: Thomson Learning, 2006.SAMUELSON, P. Introdução à Economia. New York: Mc Graw-Hill Book Company."
$new7 = "MANKIW, N.G. Introdução à economia. São Paulo: Thomson Learning, 2006.SAMUELSON, P. Introdução à Economia. New York: Mc Graw-Hill Book Company.BACHA, Edmar. Introdução à Macroeconomia: Uma perspectiva brasileira. Rio de Janeiro: Campus, 1987.BACHA et al. Estado da Economia Mundial - Desafios e Respostas - Seminário em Homenagem a Pedro Malan. São Paulo: LTC, 2015.FURTADO, C. Formação econômica do Brasil. São Paulo: Companhia Editora Nacional, 2003.GREMAUD, A. P.; VASCONCELLOS, M. A. S.; TONETO JÚNIOR, R. Economia Brasileira Contemporânea. 8 ed. São Paulo: Atlas, 2017.VASCONCELLOS, M. A. S.; GARCIA, M. E. Fundamentos de Economia. 6 ed. São Paulo: Saraiva, 2018.VASCONCELLOS, M. A. S. ECONOMIA: Micro e Macro. São Paulo: Atlas, 2015."
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2) | Out-Null
